# Saldo.xlsx update — apply the account-balance refresh described by the diff.
#
# Net changes to the "Export" sheet's data rows (Conta / Nome / Saldo):
#  - New top row: 005121919 JORGE 850978.39
#  - 004479287 ANA 43728.89            -> 004212438 KENIA 108279.35
#  - 004493324 DANIEL 7695.34  +
#    004505474 RICARDO 7551.89         -> single row 004809902 PEDRO 7376.66
#  - 004471893 PAULA 5121.13           -> 003836362 ISABELLA 3208.68
#                                          (+ new row) 005366255 RAPHAELA 2000
#  - new row inserted before 002064834 RAFAELA: 004322719 GISELA 814.09
#  - new row inserted before 005142592 ALBERTO: 004813088 JULIANA 563.32
#  - old row 004322719 GISELA 493.87   -> removed (GISELA re-appears above with new balance)
#  - old row 005121919 JORGE 297.95    -> removed (JORGE re-appears at top with new balance)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns("A")

function Find-AccountRow($account) {
    $hit = $col.Find($account)
    if ($hit -eq $null) { return -1 }
    return $hit.Row
}

function Set-TextCell($row, $colLetter, $text) {
    # Account numbers carry significant leading zeros, so they must be
    # written as text — NumberFormat "@" forces that, and ClearFormats
    # afterwards drops the now-unneeded format so the cell keeps the
    # workbook's default (unstyled) look, matching the other rows.
    $target = $ws.Range("$colLetter$row")
    $target.NumberFormat = "@"
    $target.Value = $text
    $target.ClearFormats()
}

function Set-DataRow($row, $account, $name, $balance) {
    Set-TextCell $row "A" $account
    Set-TextCell $row "B" $name
    $ws.Range("C$row").Value = $balance
}

function Insert-DataRowBefore($beforeRow, $account, $name, $balance) {
    $ws.Rows($beforeRow).Insert()
    Set-DataRow $beforeRow $account $name $balance
}

# --- 1. Remove the two rows that simply relocate elsewhere with a new balance ---
$r = Find-AccountRow "005121919"
if ($r -gt 0) { $ws.Rows($r).Delete() }

$r = Find-AccountRow "004322719"
if ($r -gt 0) { $ws.Rows($r).Delete() }

# --- 2. Update ANA -> KENIA in place ---
$r = Find-AccountRow "004479287"
Set-DataRow $r "004212438" "KENIA" 108279.35

# --- 3. Collapse DANIEL + RICARDO into a single PEDRO row ---
$r = Find-AccountRow "004493324"
Set-DataRow $r "004809902" "PEDRO" 7376.66
$r2 = Find-AccountRow "004505474"
if ($r2 -gt 0) { $ws.Rows($r2).Delete() }

# --- 4. PAULA -> ISABELLA, plus a new RAPHAELA row right after it ---
$r = Find-AccountRow "004471893"
Set-DataRow $r "003836362" "ISABELLA" 3208.68
Insert-DataRowBefore ($r + 1) "005366255" "RAPHAELA" 2000

# --- 5. New GISELA row just before RAFAELA ---
$r = Find-AccountRow "002064834"
Insert-DataRowBefore $r "004322719" "GISELA" 814.09

# --- 6. New JULIANA row just before ALBERTO ---
$r = Find-AccountRow "005142592"
Insert-DataRowBefore $r "004813088" "JULIANA" 563.32

# --- 7. New JORGE row at the very top of the data (right after the header) ---
$r = Find-AccountRow "005305448"
Insert-DataRowBefore $r "005121919" "JORGE" 850978.39
